$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.876.07"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "2.533.16"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "318.25"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").Value = "96.71"
$ws.Range("E6").Value = "  +1.01%  "

$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("D10").Value = "35.86"
$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("E11").Value = "  +0.53%  "

$ws.Range("D12").Value = "7.57"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("E13").Value = "  -3.75%  "

$ws.Range("D14").Value = "2.922.69"
$ws.Range("E14").Value = "  -0.25%  "

$ws.Range("D15").Value = "2.506.90"
$ws.Range("E15").Value = "  -1.83%  "

$ws.Range("D16").Value = "15.14"
$ws.Range("E16").Value = "  -4.16%  "

$ws.Range("D17").Value = "0.852"
$ws.Range("E17").Value = "  -1.99%  "

$ws.Range("D18").Value = "42.924.42"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("E19").Value = "  +2.91%  "

$ws.Range("D20").Value = "12.71"
$ws.Range("E20").Value = "  -3.13%  "

$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("D22").Value = "69.68"
$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("D23").Value = "253.58"
$ws.Range("E23").Value = "  -0.71%  "

$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("D25").Value = "2.06"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("D26").Value = "26.34"
$ws.Range("E26").Value = "  -4.45%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  +2.30%  "

$ws.Range("D29").Value = "40.99"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("D30").Value = "10.56"
$ws.Range("E30").Value = "  +4.34%  "

$ws.Range("E31").Value = "  -0.73%  "

$ws.Range("D32").Value = "157.73"
$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "19.39"
$ws.Range("E34").Value = "  -2.49%  "

$ws.Range("D35").Value = "2.70"
$ws.Range("E35").Value = "  +3.12%  "

$ws.Range("E36").Value = "  -1.91%  "

$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("E39").Value = "  +7.95%  "

$ws.Range("D41").Value = "21.91"
$ws.Range("E41").Value = "  -12.87%  "

$ws.Range("E42").Value = "  +0.45%  "

$ws.Range("E43").Value = "  -1.11%  "

$ws.Range("E44").Value = "  +0.21%  "

$ws.Range("E45").Value = "  -3.78%  "

$ws.Range("D46").Value = "2.009.33"
$ws.Range("E46").Value = "  -1.87%  "

$ws.Range("E47").Value = "  +2.83%  "

$ws.Range("D48").Value = "84.23"
$ws.Range("E48").Value = "  -1.95%  "

$ws.Range("D49").Value = "106.82"
$ws.Range("E49").Value = "  +4.37%  "

$ws.Range("D50").Value = "74.95"
$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("D51").Value = "2.776.18"
$ws.Range("E51").Value = "  -0.11%  "
